# Sample Project / Main.xlsx - rule row 11 ("R40") is renamed to "1".
#
# B11 must end up as a *text* cell containing "1" (not the number 1),
# while keeping its original style (s="23", General number format).
# A plain `.Value = "1"` assignment auto-coerces a numeric-looking
# string to a number, and forcing text via NumberFormat="@" (or a
# leading apostrophe) stamps the cell with a brand-new Text-formatted
# style. So: stash B11's current format on a scratch cell, force the
# text entry, then paste the stashed format back over B11 and clean
# up the scratch cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")
$scratch = $ws.Range("G11")

$target.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats - snapshot B11's current style

$target.NumberFormat = "@"
$target.Value = "1"

$scratch.Copy()
$target.PasteSpecial(-4122)    # xlPasteFormats - restore B11's original style

$scratch.Clear()
